$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.183521866798401
$ws.Range("B1").Value = 3.312803983688354
$ws.Range("C1").Value = 2.572841167449951
$ws.Range("D1").Value = 1.296128869056702
$ws.Range("E1").Value = 0.9479374289512634
